$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1796
$ws1.Range("F9").Value = 981
$ws1.Range("F14").Value = 4031
$ws1.Range("F21").Value = 181
$ws1.Range("F22").Value = 2094
$ws1.Range("F25").Value = 1977
$ws1.Range("F29").Value = 8776
$ws1.Range("F30").Value = 5798
$ws1.Range("F34").Value = 15
$ws1.Range("F35").Value = 781
$ws1.Range("F37").Value = 49
$ws1.Range("F42").Value = 174
$ws1.Range("F43").Value = 4667
$ws1.Range("F47").Value = 416

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8409

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8409
$ws4.Range("F7").Value = 1796
$ws4.Range("F12").Value = 981
$ws4.Range("F15").Value = 4031
$ws4.Range("F21").Value = 181
$ws4.Range("F22").Value = 2094
$ws4.Range("F28").Value = 1977
$ws4.Range("F33").Value = 8776
$ws4.Range("F36").Value = 781
$ws4.Range("F42").Value = 174
$ws4.Range("F45").Value = 416

Write-Output "Updated 26 cells across sheets 展览, 本地生活, 全部类型"
